$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as exact text, avoiding Excel auto-number/date coercion,
# and clear the temporary text-format style afterwards so no stray style id is left on the cell.
function Set-ExactText {
    param($range, $text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-ExactText $ws.Range('D2') '26.930.99'
Set-ExactText $ws.Range('E2') '  -1.67%  '
Set-ExactText $ws.Range('D3') '1.813.48'
Set-ExactText $ws.Range('E3') '  -0.68%  '
Set-ExactText $ws.Range('D4') '0.9980'
Set-ExactText $ws.Range('E4') '  -0.28%  '
Set-ExactText $ws.Range('D5') '309.91'
Set-ExactText $ws.Range('E5') '  -1.34%  '
Set-ExactText $ws.Range('D6') '0.9983'
Set-ExactText $ws.Range('E6') '  -0.29%  '
Set-ExactText $ws.Range('D7') '0.4615'
Set-ExactText $ws.Range('D8') '0.3747'
Set-ExactText $ws.Range('E8') '  -0.33%  '
Set-ExactText $ws.Range('D9') '0.07325'
Set-ExactText $ws.Range('E9') '  -2.53%  '
Set-ExactText $ws.Range('D10') '0.8694'
Set-ExactText $ws.Range('E10') '  -2.42%  '
Set-ExactText $ws.Range('D11') '20.47'
Set-ExactText $ws.Range('E11') '  -2.97%  '
Set-ExactText $ws.Range('D12') '1.745.79'
Set-ExactText $ws.Range('E12') '  -4.43%  '
Set-ExactText $ws.Range('D13') '5.343'
Set-ExactText $ws.Range('E13') '  -1.35%  '
Set-ExactText $ws.Range('D14') '6.529'
Set-ExactText $ws.Range('E14') '  -3.46%  '
Set-ExactText $ws.Range('D15') '0.07040'
Set-ExactText $ws.Range('E15') '  -0.99%  '
Set-ExactText $ws.Range('E16') '  -2.93%  '
Set-ExactText $ws.Range('D17') '0.9989'
Set-ExactText $ws.Range('E17') '  -0.34%  '
Set-ExactText $ws.Range('D18') '0.000008696'
Set-ExactText $ws.Range('E18') '  -1.28%  '
Set-ExactText $ws.Range('D19') '0.9977'
Set-ExactText $ws.Range('E19') '  -0.25%  '
Set-ExactText $ws.Range('D20') '14.77'
Set-ExactText $ws.Range('E20') '  -2.94%  '
Set-ExactText $ws.Range('D21') '26.917.44'
Set-ExactText $ws.Range('E21') '  -1.71%  '
Set-ExactText $ws.Range('D22') '5.302'
Set-ExactText $ws.Range('E22') '  +0.71%  '
Set-ExactText $ws.Range('D23') '10.71'
Set-ExactText $ws.Range('D24') '1.967.01'
Set-ExactText $ws.Range('E24') '  -4.42%  '
Set-ExactText $ws.Range('D25') '1.911'
Set-ExactText $ws.Range('E25') '  -3.27%  '
Set-ExactText $ws.Range('D26') '150.86'
Set-ExactText $ws.Range('E26') '  -0.43%  '
Set-ExactText $ws.Range('D27') '18.37'
Set-ExactText $ws.Range('E27') '  -1.20%  '
Set-ExactText $ws.Range('D28') '2.160'
Set-ExactText $ws.Range('E28') '  -9.27%  '
Set-ExactText $ws.Range('D29') '5.269'
Set-ExactText $ws.Range('E29') '  -1.79%  '
Set-ExactText $ws.Range('D30') '114.92'
Set-ExactText $ws.Range('E30') '  -2.39%  '
Set-ExactText $ws.Range('D31') '0.08888'
Set-ExactText $ws.Range('E31') '  +0.67%  '
Set-ExactText $ws.Range('D32') '0.7670'
Set-ExactText $ws.Range('E32') '  -2.61%  '
Set-ExactText $ws.Range('D33') '1.172'
Set-ExactText $ws.Range('E33') '  -2.48%  '
Set-ExactText $ws.Range('D34') '4.473'
Set-ExactText $ws.Range('E34') '  -1.11%  '
Set-ExactText $ws.Range('D35') '2.887'
Set-ExactText $ws.Range('E35') '  -0.59%  '
Set-ExactText $ws.Range('D36') '0.9975'
Set-ExactText $ws.Range('E36') '  -0.36%  '
Set-ExactText $ws.Range('D37') '1.119'
Set-ExactText $ws.Range('E37') '  +0.74%  '
Set-ExactText $ws.Range('B38') 'RenderToken'
Set-ExactText $ws.Range('C38') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-ExactText $ws.Range('D38') '2.495'
Set-ExactText $ws.Range('E38') '  +9.15%  '
Set-ExactText $ws.Range('B39') 'VeChain'
Set-ExactText $ws.Range('C39') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-ExactText $ws.Range('D39') '0.01951'
Set-ExactText $ws.Range('E39') '  -2.11%  '
Set-ExactText $ws.Range('D40') '0.05234'
Set-ExactText $ws.Range('E40') '  -1.96%  '
Set-ExactText $ws.Range('D41') '2.907'
Set-ExactText $ws.Range('E41') '  +1.60%  '
Set-ExactText $ws.Range('D42') '7.182'
Set-ExactText $ws.Range('E42') '  -2.84%  '
Set-ExactText $ws.Range('D43') '0.5262'
Set-ExactText $ws.Range('E43') '  -1.08%  '
Set-ExactText $ws.Range('D44') '0.1659'
Set-ExactText $ws.Range('E44') '  -4.26%  '
Set-ExactText $ws.Range('D45') '8.568'
Set-ExactText $ws.Range('E45') '  -2.18%  '
Set-ExactText $ws.Range('D46') '0.5052'
Set-ExactText $ws.Range('E46') '  -1.23%  '
Set-ExactText $ws.Range('D47') '10.28'
Set-ExactText $ws.Range('E47') '  -3.52%  '
Set-ExactText $ws.Range('D48') '104.35'
Set-ExactText $ws.Range('E48') '  -1.75%  '
Set-ExactText $ws.Range('D49') '0.9976'
Set-ExactText $ws.Range('E49') '  -0.33%  '
Set-ExactText $ws.Range('D50') '1.667'
Set-ExactText $ws.Range('E50') '  -2.25%  '
Set-ExactText $ws.Range('D51') '0.06317'
Set-ExactText $ws.Range('E51') '  -0.94%  '
